$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.963.45"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "'3.791.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.80%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'616.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.40%  "
$ws.Range("E6").Value = "  -3.81%  "
$ws.Range("D7").Value = "'3.789.49"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.89%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +0.45%  "
$ws.Range("D10").Value = "'0.168"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.79%  "
$ws.Range("D11").Value = "'6.34"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.09%  "
$ws.Range("D12").Value = "'0.494"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("D13").Value = "'41.22"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.09%  "
$ws.Range("E14").Value = "  +1.38%  "
$ws.Range("D15").Value = "'4.422.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.60%  "
$ws.Range("D16").Value = "'3.783.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.53%  "
$ws.Range("D17").Value = "'70.042.02"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.19%  "
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").Value = "'7.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.84%  "
$ws.Range("D20").Value = "'515.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.19%  "
$ws.Range("D21").Value = "'16.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.86%  "
$ws.Range("E22").Value = "  +3.02%  "
$ws.Range("D23").Value = "'0.729"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.23%  "
$ws.Range("D24").Value = "'2.53"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.27%  "
$ws.Range("D25").Value = "'88.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.30%  "
$ws.Range("D26").Value = "'13.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.13%  "
$ws.Range("D27").Value = "'11.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.61%  "
$ws.Range("D28").Value = "'0.0000134"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +24.91%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  -1.61%  "
$ws.Range("E31").Value = "  +3.27%  "
$ws.Range("D32").Value = "'7.82"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.90%  "
$ws.Range("D33").Value = "'31.77"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.21%  "
$ws.Range("D34").Value = "'0.116"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.60%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.24%  "
$ws.Range("D36").Value = "'6.24"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.28%  "
$ws.Range("E37").Value = "  +2.11%  "
$ws.Range("D38").Value = "'0.341"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.79%  "
$ws.Range("E39").Value = "  +3.43%  "
$ws.Range("E40").Value = "  +3.60%  "
$ws.Range("D41").Value = "'51.39"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.00%  "
$ws.Range("D42").Value = "'44.47"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.51%  "
$ws.Range("E43").Value = "  -0.36%  "
$ws.Range("D44").Value = "'424.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.93%  "
$ws.Range("D45").Value = "'3.073.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.90%  "
$ws.Range("E46").Value = "  -0.86%  "
$ws.Range("D47").Value = "'0.0365"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.42%  "
$ws.Range("D48").Value = "'27.74"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("D49").Value = "'2.51"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.22%  "
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("D51").Value = "'135.44"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.93%  "
